# cierre 18 Jun 22
$wb = $excel.ActiveWorkbook

$wsArq = $wb.Worksheets.Item("ARQUITECTO        ")

# Correct the written-out amount: "CIENTO MIL" (one hundred thousand) -> "CIEN MIL" (one hundred thousand, short form)
$wsArq.Range("A2").Value = "CIEN    MIL   PESOS 00/100 M.N."

# Bonus amount updated from 50000 to 100000
$wsArq.Range("D1").Value = 100000

# Move/leave the on-screen selection on A11:B11 (active cell A11) for the closing snapshot
$wsArq.Activate() | Out-Null
$wsArq.Range("A11:B11").Select() | Out-Null

# Recalculate so the TODAY() cells (A11 on both vale sheets) pick up the closing date
$excel.CalculateFull()
